$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 210.07692  # H28
$ws.Cells.Item(28, 9).Value = 148.27272  # I28
$ws.Cells.Item(28, 10).Value = 550  # J28
$ws.Cells.Item(28, 11).Value = 148.27272  # K28
$ws.Cells.Item(28, 12).Value = 550  # L28
$ws.Cells.Item(28, 13).Value = 336.72728  # M28
$ws.Cells.Item(28, 14).Value = -1520  # N28

$ws.Cells.Item(39, 8).Value = 696.2778  # H39
$ws.Cells.Item(39, 9).Value = 94.818184  # I39
$ws.Cells.Item(39, 10).Value = 1641.4286  # J39
$ws.Cells.Item(39, 11).Value = 284.454552  # K39
$ws.Cells.Item(39, 12).Value = 4924.2858  # L39
$ws.Cells.Item(39, 13).Value = 11.54544799999996  # M39
$ws.Cells.Item(39, 14).Value = -5516.2858  # N39

$ws.Cells.Item(43, 8).Value = 800  # H43
$ws.Cells.Item(43, 9).Value = 600  # I43
$ws.Cells.Item(43, 10).Value = 1000  # J43
$ws.Cells.Item(43, 11).Value = 600  # K43
$ws.Cells.Item(43, 12).Value = 1000  # L43
$ws.Cells.Item(43, 13).Value = -531  # M43
$ws.Cells.Item(43, 14).Value = -1138  # N43

$ws.Cells.Item(70, 8).Value = 34050  # H70
$ws.Cells.Item(70, 9).Value = 1200  # I70
$ws.Cells.Item(70, 10).Value = 45000  # J70
$ws.Cells.Item(70, 11).Value = 3600  # K70
$ws.Cells.Item(70, 12).Value = 135000  # L70
$ws.Cells.Item(70, 13).Value = -3330  # M70
$ws.Cells.Item(70, 14).Value = -135540  # N70

$ws.Cells.Item(73, 8).Value = 34050  # H73
$ws.Cells.Item(73, 9).Value = 1200  # I73
$ws.Cells.Item(73, 10).Value = 45000  # J73
$ws.Cells.Item(73, 11).Value = 3600  # K73
$ws.Cells.Item(73, 12).Value = 135000  # L73
$ws.Cells.Item(73, 13).Value = -2664  # M73
$ws.Cells.Item(73, 14).Value = -136872  # N73

$ws.Cells.Item(74, 8).Value = 3498.8572  # H74
$ws.Cells.Item(74, 9).Value = 2373.25  # I74
$ws.Cells.Item(74, 11).Value = 2373.25  # K74
$ws.Cells.Item(74, 13).Value = -1437.25  # M74

$ws.Cells.Item(77, 8).Value = 3498.8572  # H77
$ws.Cells.Item(77, 9).Value = 2373.25  # I77
$ws.Cells.Item(77, 11).Value = 11866.25  # K77
$ws.Cells.Item(77, 13).Value = -7186.25  # M77

$ws.Cells.Item(125, 8).Value = 1140  # H125
$ws.Cells.Item(125, 9).Value = 1140  # I125
$ws.Cells.Item(125, 11).Value = 10260  # K125
$ws.Cells.Item(125, 13).Value = -7800  # M125

$ws.Cells.Item(135, 8).Value = 2034.8  # H135
$ws.Cells.Item(135, 9).Value = 2033  # I135
$ws.Cells.Item(135, 11).Value = 18297  # K135
$ws.Cells.Item(135, 13).Value = -15762  # M135

$ws.Cells.Item(138, 8).Value = 1849.463  # H138
$ws.Cells.Item(138, 9).Value = 1250.7307  # I138
$ws.Cells.Item(138, 10).Value = 2405.4285  # J138
$ws.Cells.Item(138, 11).Value = 3752.1921  # K138
$ws.Cells.Item(138, 12).Value = 7216.2855  # L138
$ws.Cells.Item(138, 13).Value = 1387.8079  # M138
$ws.Cells.Item(138, 14).Value = -17496.2855  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 327851.53  # H2
$ws.Cells.Item(2, 9).Value = 464017.5  # I2
$ws.Cells.Item(2, 10).Value = 1053.2  # J2
$ws.Cells.Item(2, 11).Value = 464017.5  # K2
$ws.Cells.Item(2, 12).Value = 1053.2  # L2
$ws.Cells.Item(2, 13).Value = -463904.5  # M2
$ws.Cells.Item(2, 14).Value = -1279.2  # N2

$ws.Cells.Item(61, 8).Value = 5626.4614  # H61
$ws.Cells.Item(61, 9).Value = 5894  # I61
$ws.Cells.Item(61, 10).Value = 5024.5  # J61
$ws.Cells.Item(61, 11).Value = 5894  # K61
$ws.Cells.Item(61, 12).Value = 5024.5  # L61
$ws.Cells.Item(61, 13).Value = -5682  # M61
$ws.Cells.Item(61, 14).Value = -5448.5  # N61

$ws.Cells.Item(74, 8).Value = 1157.5  # H74
$ws.Cells.Item(74, 9).Value = 559  # I74
$ws.Cells.Item(74, 11).Value = 559  # K74
$ws.Cells.Item(74, 13).Value = 315  # M74

$ws.Cells.Item(77, 8).Value = 1157.5  # H77
$ws.Cells.Item(77, 9).Value = 559  # I77
$ws.Cells.Item(77, 11).Value = 2795  # K77
$ws.Cells.Item(77, 13).Value = 1573  # M77

$ws.Cells.Item(97, 8).Value = 962.56525  # H97
$ws.Cells.Item(97, 9).Value = 933.63635  # I97
$ws.Cells.Item(97, 11).Value = 933.63635  # K97
$ws.Cells.Item(97, 13).Value = -437.63635  # M97

$ws.Cells.Item(116, 8).Value = 327851.53  # H116
$ws.Cells.Item(116, 9).Value = 464017.5  # I116
$ws.Cells.Item(116, 10).Value = 1053.2  # J116
$ws.Cells.Item(116, 11).Value = 464017.5  # K116
$ws.Cells.Item(116, 12).Value = 1053.2  # L116
$ws.Cells.Item(116, 13).Value = -461723.5  # M116
$ws.Cells.Item(116, 14).Value = -5641.2  # N116

$ws.Cells.Item(132, 8).Value = 1777.4348  # H132
$ws.Cells.Item(132, 9).Value = 1257.0588  # I132
$ws.Cells.Item(132, 10).Value = 3251.8333  # J132
$ws.Cells.Item(132, 11).Value = 3771.1764  # K132
$ws.Cells.Item(132, 12).Value = 9755.499899999999  # L132
$ws.Cells.Item(132, 13).Value = -1241.1764  # M132
$ws.Cells.Item(132, 14).Value = -14815.4999  # N132

$ws.Cells.Item(136, 8).Value = 5626.4614  # H136
$ws.Cells.Item(136, 9).Value = 5894  # I136
$ws.Cells.Item(136, 10).Value = 5024.5  # J136
$ws.Cells.Item(136, 11).Value = 17682  # K136
$ws.Cells.Item(136, 12).Value = 15073.5  # L136
$ws.Cells.Item(136, 13).Value = -15132  # M136
$ws.Cells.Item(136, 14).Value = -20173.5  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 327851.53  # H3
$ws.Cells.Item(3, 9).Value = 464017.5  # I3
$ws.Cells.Item(3, 10).Value = 1053.2  # J3
$ws.Cells.Item(3, 11).Value = 464017.5  # K3
$ws.Cells.Item(3, 12).Value = 1053.2  # L3
$ws.Cells.Item(3, 13).Value = -463903.5  # M3
$ws.Cells.Item(3, 14).Value = -1281.2  # N3

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1616.3334  # H22
$ws.Cells.Item(22, 9).Value = 500  # I22
$ws.Cells.Item(22, 10).Value = 1717.8182  # J22
$ws.Cells.Item(22, 11).Value = 500  # K22
$ws.Cells.Item(22, 12).Value = 1717.8182  # L22
$ws.Cells.Item(22, 13).Value = -150  # M22
$ws.Cells.Item(22, 14).Value = -2417.8182  # N22

$ws.Cells.Item(31, 8).Value = 2321.8333  # H31
$ws.Cells.Item(31, 9).Value = 2178.6  # I31
$ws.Cells.Item(31, 11).Value = 2178.6  # K31
$ws.Cells.Item(31, 13).Value = -1883.6  # M31

$ws.Cells.Item(34, 8).Value = 2321.8333  # H34
$ws.Cells.Item(34, 9).Value = 2178.6  # I34
$ws.Cells.Item(34, 11).Value = 2178.6  # K34
$ws.Cells.Item(34, 13).Value = -1976.6  # M34

$ws.Cells.Item(58, 8).Value = 3625355.8  # H58
$ws.Cells.Item(58, 9).Value = 6212009.5  # I58
$ws.Cells.Item(58, 11).Value = 6212009.5  # K58
$ws.Cells.Item(58, 13).Value = -6211806.5  # M58

$ws.Cells.Item(62, 8).Value = 2965.3333  # H62
$ws.Cells.Item(62, 9).Value = 2948.25  # I62
$ws.Cells.Item(62, 10).Value = 2999.5  # J62
$ws.Cells.Item(62, 11).Value = 2948.25  # K62
$ws.Cells.Item(62, 12).Value = 2999.5  # L62
$ws.Cells.Item(62, 13).Value = -2324.25  # M62
$ws.Cells.Item(62, 14).Value = -4247.5  # N62

$ws.Cells.Item(65, 8).Value = 2965.3333  # H65
$ws.Cells.Item(65, 9).Value = 2948.25  # I65
$ws.Cells.Item(65, 10).Value = 2999.5  # J65
$ws.Cells.Item(65, 11).Value = 14741.25  # K65
$ws.Cells.Item(65, 12).Value = 14997.5  # L65
$ws.Cells.Item(65, 13).Value = -11621.25  # M65
$ws.Cells.Item(65, 14).Value = -21237.5  # N65

$ws.Cells.Item(107, 8).Value = 1182  # H107
$ws.Cells.Item(107, 9).Value = 859.75  # I107
$ws.Cells.Item(107, 10).Value = 1439.8  # J107
$ws.Cells.Item(107, 11).Value = 859.75  # K107
$ws.Cells.Item(107, 12).Value = 1439.8  # L107
$ws.Cells.Item(107, 13).Value = 1060.25  # M107
$ws.Cells.Item(107, 14).Value = -5279.8  # N107

$ws.Cells.Item(134, 8).Value = 2938.5833  # H134
$ws.Cells.Item(134, 9).Value = 2381.889  # I134
$ws.Cells.Item(134, 11).Value = 7145.667  # K134
$ws.Cells.Item(134, 13).Value = -4610.667  # M134

$ws.Cells.Item(136, 8).Value = 3625355.8  # H136
$ws.Cells.Item(136, 9).Value = 6212009.5  # I136
$ws.Cells.Item(136, 11).Value = 18636028.5  # K136
$ws.Cells.Item(136, 13).Value = -18633478.5  # M136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 274.75  # H2
$ws.Cells.Item(2, 9).Value = 191.66667  # I2
$ws.Cells.Item(2, 10).Value = 524  # J2
$ws.Cells.Item(2, 11).Value = 1150.00002  # K2
$ws.Cells.Item(2, 12).Value = 3144  # L2
$ws.Cells.Item(2, 13).Value = -1037.00002  # M2
$ws.Cells.Item(2, 14).Value = -3370  # N2

$ws.Cells.Item(7, 8).Value = 631.63635  # H7
$ws.Cells.Item(7, 10).Value = 705.3333  # J7
$ws.Cells.Item(7, 12).Value = 2115.9999  # L7
$ws.Cells.Item(7, 14).Value = -2339.9999  # N7

$ws.Cells.Item(38, 8).Value = 552.4545000000001  # H38
$ws.Cells.Item(38, 9).Value = 153.71428  # I38
$ws.Cells.Item(38, 11).Value = 461.14284  # K38
$ws.Cells.Item(38, 13).Value = -114.14284  # M38

$ws.Cells.Item(131, 8).Value = 19648.514  # H131
$ws.Cells.Item(131, 10).Value = 21955  # J131
$ws.Cells.Item(131, 12).Value = 65865  # L131
$ws.Cells.Item(131, 14).Value = -75945  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2000  # H80
$ws.Cells.Item(80, 10).Value = 2000  # J80
$ws.Cells.Item(80, 12).Value = 2000  # L80
$ws.Cells.Item(80, 14).Value = -3996  # N80

$ws.Cells.Item(83, 8).Value = 2000  # H83
$ws.Cells.Item(83, 10).Value = 2000  # J83
$ws.Cells.Item(83, 12).Value = 10000  # L83
$ws.Cells.Item(83, 14).Value = -19984  # N83

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 12978.8  # H16
$ws.Cells.Item(16, 9).Value = 12978.8  # I16
$ws.Cells.Item(16, 10).Value = 0  # J16
$ws.Cells.Item(16, 11).Value = 12978.8  # K16
$ws.Cells.Item(16, 12).Value = 0  # L16
$ws.Cells.Item(16, 13).Value = -12808.8  # M16
$ws.Cells.Item(16, 14).ClearContents()  # N16 removed

$ws.Cells.Item(96, 8).Value = 75000  # H96
$ws.Cells.Item(96, 10).Value = 75000  # J96
$ws.Cells.Item(96, 12).Value = 75000  # L96
$ws.Cells.Item(96, 14).Value = -80492  # N96

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 19842614  # H136
$ws.Cells.Item(136, 9).Value = 27778746  # I136
$ws.Cells.Item(136, 11).Value = 83336238  # K136
$ws.Cells.Item(136, 13).Value = -83333688  # M136
